$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "runs" (C) and "balls" (D) figures for Kartik Tyagi's 2nd and 4th
# score lines were swapped. Row 2 should carry the 2/3 that was
# mistakenly left on row 4, and row 4 should go back to 0/0.
#
# These columns store numbers-as-text (the sheet already carries a
# numberStoredAsText ignored-error hint for A1:F4), so force each cell
# back to Text before writing the new value and then drop the
# temporary number format again so the cell's style is left untouched.

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("C2") "2"
Set-TextValue $ws.Range("D2") "3"
Set-TextValue $ws.Range("C4") "0"
Set-TextValue $ws.Range("D4") "0"
